$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width changes ---
# Note: Excel's ColumnWidth property applies an offset of +5/6 (0.8333...)
# relative to the raw OOXML "width" attribute value when saved, so we
# subtract that offset here to land on the exact target stored widths.
# col min=1 max=1: width 25 -> 28  (column A)
$ws.Columns.Item(1).ColumnWidth = 27.16666666666667
# NEW col min=10 max=10: width 25  (column J)
$ws.Columns.Item(10).ColumnWidth = 24.16666666666667
# col min=12 max=12: width 30 -> 25  (column L)
$ws.Columns.Item(12).ColumnWidth = 24.16666666666667

# --- I1 cell text change ---
$ws.Range("I1").Value = "MEDE-CABA-H-07-COS"

# --- L4, L5, L7, L8 cells get the new text value ---
$ws.Range("L4").Value = "MEDE-CABA-H-07-COS"
$ws.Range("L5").Value = "MEDE-CABA-H-07-COS"
$ws.Range("L7").Value = "MEDE-CABA-H-07-COS"
$ws.Range("L8").Value = "MEDE-CABA-H-07-COS"

# --- F20-F34 interface/description/logical-channel text updates ---
$ws.Range("F20").Value = "interface upstream 0/18.0"
$ws.Range("F21").Value = '  description "PUERTO LIBRE"'
$ws.Range("F22").Value = '  logical-channel 0 description "PUERTO LIBRE"'

$ws.Range("F24").Value = "interface upstream 0/18.1"
$ws.Range("F25").Value = '  description "PUERTO LIBRE"'
$ws.Range("F26").Value = '  logical-channel 0 description "PUERTO LIBRE"'

$ws.Range("F28").Value = "interface upstream 0/18.2"
$ws.Range("F29").Value = '  description "PUERTO LIBRE"'
$ws.Range("F30").Value = '  logical-channel 0 description "PUERTO LIBRE"'

$ws.Range("F32").Value = "interface upstream 0/18.3"
$ws.Range("F33").Value = '  description "PUERTO LIBRE"'
$ws.Range("F34").Value = '  logical-channel 0 description "PUERTO LIBRE"'

# --- F39 text update ---
$ws.Range("F39").Value = "no service group  IRL "

# --- Merge cells change: G4:G8 -> G4:G7 ---
$ws.Range("G4:G8").UnMerge()
$ws.Range("G4:G7").Merge()
